$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.24"
$ws.Range("E2").Value = "'-2.30%"

$ws.Range("E3").Value = "'-1.68%"

$ws.Range("D4").Value = "'5.143"
$ws.Range("E4").Value = "'-2.58%"

$ws.Range("E5").Value = "'-3.05%"

$ws.Range("D6").Value = "'1.843"
$ws.Range("E6").Value = "'24.25%"

$ws.Range("E7").Value = "'-1.25%"

$ws.Range("E8").Value = "'-0.42%"

$ws.Range("D9").Value = "'0.9257"
$ws.Range("E9").Value = "'1.09%"

$ws.Range("D10").Value = "'0.1678"
$ws.Range("E10").Value = "'-0.95%"

$ws.Range("D11").Value = "'0.07108"
$ws.Range("E11").Value = "'-8.47%"

$ws.Range("D12").Value = "'0.08007"
$ws.Range("E12").Value = "'-0.96%"

$ws.Range("D13").Value = "'0.03000"
$ws.Range("E13").Value = "'-0.52%"

$ws.Range("D14").Value = "'0.09923"
$ws.Range("E14").Value = "'0.35%"

$ws.Range("D15").Value = "'0.001498"
$ws.Range("E15").Value = "'0.69%"

$ws.Range("D16").Value = "'0.006190"
$ws.Range("E16").Value = "'0.65%"

$ws.Range("E17").Value = "'-0.69%"

$ws.Range("D18").Value = "'2.222"
$ws.Range("E18").Value = "'-0.44%"

$ws.Range("E19").Value = "'-2.29%"

$ws.Range("E20").Value = "'-1.37%"

$ws.Range("D21").Value = "'4.563"
$ws.Range("E21").Value = "'1.69%"

$ws.Range("D22").Value = "'0.04642"
$ws.Range("E22").Value = "'2.13%"

$ws.Range("E23").Value = "'-3.22%"

$ws.Range("E24").Value = "'0.12%"

$ws.Range("D25").Value = "'0.004729"
$ws.Range("E25").Value = "'6.57%"

$ws.Range("E26").Value = "'-7.00%"

$ws.Range("D27").Value = "'0.0001875"
$ws.Range("E27").Value = "'7.91%"

$ws.Range("D39").Value = "'0.01711"
$ws.Range("E39").Value = "'-3.52%"

$ws.Range("D40").Value = "'0.04472"
$ws.Range("E40").Value = "'-1.29%"

$ws.Range("D41").Value = "'0.007151"
$ws.Range("E41").Value = "'-1.06%"

$ws.Range("D42").Value = "'0.1332"
$ws.Range("E42").Value = "'-1.00%"

$ws.Range("D43").Value = "'0.002149"
$ws.Range("E43").Value = "'-3.87%"

$ws.Range("D44").Value = "'0.01048"
$ws.Range("E44").Value = "'-22.29%"

$ws.Range("D45").Value = "'0.00006222"
$ws.Range("E45").Value = "'0.02%"

$ws.Range("E46").Value = "'-21.26%"

$ws.Range("D47").Value = "'0.7393"
$ws.Range("E47").Value = "'4.25%"
